# Generate Report for handoff
#
# The old e2e markdown file (6cf34fa2-...) was replaced by a new one
# (bf5cc28a-...), and a second source file (5b906697-...) was added whose
# handoff transform failed. This pushes the ".localization-config" /
# "Not to be localized" row down by one on every sheet, and refreshes the
# handoff file name / timestamp for the file that succeeded.

$wb = $excel.ActiveWorkbook

$oldMd   = "6cf34fa2-e9c1-4580-9a61-c6ccfd397e1e.md"
$newMd   = "bf5cc28a-226f-4ef1-81c1-6de60694074d.md"
$failMd  = "5b906697-7601-481e-ac4a-b63580308d43.md"
$cfgName = ".localization-config"

$newXlfZh = "bf5cc28a-226f-4ef1-81c1-6de60694074d.307a61312796808768d36c0754601fdec402bd9f.zh-cn.xlf"
$newXlfDe = "bf5cc28a-226f-4ef1-81c1-6de60694074d.307a61312796808768d36c0754601fdec402bd9f.de-de.xlf"

$newZhTime = "2016-01-25 13:58:27"
$newDeTime = "2016-01-25 13:58:36"
$epoch     = "0001-01-01 00:00:00"
$dateFmt   = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value2 = $newMd
$ws.Range("B2").Value2 = "Ready for handoff"
$ws.Range("C2").Value2 = "Ready for handoff"

$ws.Range("A3").Value2 = $failMd
$ws.Range("B3").Value2 = "Handoff transform failed"
$ws.Range("C3").Value2 = "Handoff transform failed"

$ws.Range("A4").Value2 = $cfgName
$ws.Range("B4").Value2 = "Not to be localized"
$ws.Range("C4").Value2 = "Not to be localized"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0896e29b171f2204e9b641812ef0078c6767e5a9/e2e/$newMd", "", "", $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0896e29b171f2204e9b641812ef0078c6767e5a9/e2e/$failMd", "", "", $failMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0896e29b171f2204e9b641812ef0078c6767e5a9/$cfgName", "", "", $cfgName) | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value2 = $newMd
$ws.Range("B2").Value2 = "Ready for handoff"
$ws.Range("C2").Value2 = $newXlfZh
$ws.Range("D2").Value2 = $newZhTime
$ws.Range("G2").Value2 = $epoch
$ws.Range("H2").Value2 = "Include"

$ws.Range("A3").Value2 = $failMd
$ws.Range("B3").Value2 = "Handoff transform failed"
$ws.Range("D3").Value2 = $epoch
$ws.Range("G3").Value2 = $epoch
$ws.Range("H3").Value2 = "Ignored"

$ws.Range("A4").Value2 = $cfgName
$ws.Range("B4").Value2 = "Not to be localized"
$ws.Range("D4").Value2 = $epoch
$ws.Range("D4").NumberFormat = $dateFmt
$ws.Range("G4").Value2 = $epoch
$ws.Range("H4").Value2 = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0896e29b171f2204e9b641812ef0078c6767e5a9/e2e/$newMd", "", "", $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7d3aac863244ae2ae02b0ae492c271637486fb9e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$newXlfZh", "", "", $newXlfZh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0896e29b171f2204e9b641812ef0078c6767e5a9/e2e/$failMd", "", "", $failMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0896e29b171f2204e9b641812ef0078c6767e5a9/$cfgName", "", "", $cfgName) | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value2 = $newMd
$ws.Range("B2").Value2 = "Ready for handoff"
$ws.Range("C2").Value2 = $newXlfDe
$ws.Range("D2").Value2 = $newDeTime
$ws.Range("G2").Value2 = $epoch
$ws.Range("H2").Value2 = "Include"

$ws.Range("A3").Value2 = $failMd
$ws.Range("B3").Value2 = "Handoff transform failed"
$ws.Range("D3").Value2 = $epoch
$ws.Range("G3").Value2 = $epoch
$ws.Range("H3").Value2 = "Ignored"

$ws.Range("A4").Value2 = $cfgName
$ws.Range("B4").Value2 = "Not to be localized"
$ws.Range("D4").Value2 = $epoch
$ws.Range("D4").NumberFormat = $dateFmt
$ws.Range("G4").Value2 = $epoch
$ws.Range("H4").Value2 = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0896e29b171f2204e9b641812ef0078c6767e5a9/e2e/$newMd", "", "", $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ab6418c4f68255506749adf541136dbe04de1474/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$newXlfDe", "", "", $newXlfDe) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0896e29b171f2204e9b641812ef0078c6767e5a9/e2e/$failMd", "", "", $failMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0896e29b171f2204e9b641812ef0078c6767e5a9/$cfgName", "", "", $cfgName) | Out-Null
